$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 83
$ws.Range("A83").Value = 11
$ws.Range("B83").Value = "Vega Monumental Concepción"
$ws.Range("C83").Value = "Bíobío"
$ws.Range("D83").Value = 44911
$ws.Range("D83").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E83").Value = 8
$ws.Range("F83").Value = 100112037
$ws.Range("G83").Value = "Cebollín"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 200
$ws.Range("K83").Value = 700
$ws.Range("L83").Value = 800
$ws.Range("M83").Value = 750
$ws.Range("N83").Value = "$/paquete 6 unidades"
$ws.Range("O83").Value = "Región de Ñuble"
$ws.Range("P83").Value = 125
$ws.Range("Q83").Value = 6
$ws.Range("R83").Value = "Hortaliza"

# Row 84
$ws.Range("A84").Value = 11
$ws.Range("B84").Value = "Vega Monumental Concepción"
$ws.Range("C84").Value = "Bíobío"
$ws.Range("D84").Value = 44911
$ws.Range("D84").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E84").Value = 8
$ws.Range("F84").Value = 100112037
$ws.Range("G84").Value = "Cebollín"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Segunda"
$ws.Range("J84").Value = 100
$ws.Range("K84").Value = 600
$ws.Range("L84").Value = 600
$ws.Range("M84").Value = 600
$ws.Range("N84").Value = "$/paquete 6 unidades"
$ws.Range("O84").Value = "Región de Ñuble"
$ws.Range("P84").Value = 100
$ws.Range("Q84").Value = 6
$ws.Range("R84").Value = "Hortaliza"
